$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 487, shifting the existing rows (487:498) down to (490:501)
$ws.Rows("487:489").Insert()

# New weekly data rows (date 2021-09-09 = serial 44448), same market/product block
$newRows = @(
    @{ Row = 487; L = "Especial"; M = 4;  N = 210000; O = 210000; P = 210000; Q = "$/bins (450 kilos)"; R = "Región de O'Higgins"; S = 467; T = 450 },
    @{ Row = 488; L = "Primera";  M = 6;  N = 190000; O = 190000; P = 190000; Q = "$/bins (450 kilos)"; R = "Región de O'Higgins"; S = 422; T = 450 },
    @{ Row = 489; L = "Segunda";  M = 8;  N = 160000; O = 160000; P = 160000; Q = "$/bins (450 kilos)"; R = "Región de O'Higgins"; S = 356; T = 450 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 9
    $ws.Cells.Item($row, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($row, 3).Value = "Metropolitana"
    $ws.Cells.Item($row, 4).Value = 44448
    $ws.Cells.Item($row, 5).Value = 13
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100104
    $ws.Cells.Item($row, 8).Value = "Frutos de pepita"
    $ws.Cells.Item($row, 9).Value = 100104005
    $ws.Cells.Item($row, 10).Value = "Pera"
    $ws.Cells.Item($row, 11).Value = "Packham's Triumph"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
